$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived NATMI metrics (row, column, new value)
$updates = @(
    @(2, 7, 29.22757333333334),
    @(2, 8, 87.68272),
    @(2, 9, 0.08948272176993048),
    @(2, 10, 0.08948272176993047),
    @(2, 13, 15.47987166666667),
    @(2, 14, 46.439615),
    @(2, 15, 0.7960757698994193),
    @(2, 16, 0.7960757698994194),
    @(2, 17, 452.4390843280889),
    @(2, 18, 4071.9517589528),
    @(2, 19, 0.07123502662569293),
    @(2, 20, 0.07123502662569293),
    @(3, 7, 29.22757333333334),
    @(3, 8, 87.68272),
    @(3, 9, 0.08948272176993048),
    @(3, 10, 0.08948272176993047),
    @(3, 15, 0.05100527512565552),
    @(3, 16, 0.05100527512565553),
    @(3, 17, 28.98817028016),
    @(3, 18, 260.89353252144),
    @(3, 19, 0.004564090842867789),
    @(3, 20, 0.004564090842867789),
    @(4, 7, 29.22757333333334),
    @(4, 8, 87.68272),
    @(4, 9, 0.08948272176993048),
    @(4, 10, 0.08948272176993047),
    @(4, 13, 0.171678),
    @(4, 14, 0.515034),
    @(4, 15, 0.008828800326496623),
    @(4, 16, 0.008828800326496624),
    @(4, 17, 5.017731334720001),
    @(4, 18, 45.15958201248),
    @(4, 19, 0.0007900250831781686),
    @(4, 20, 0.0007900250831781688),
    @(5, 7, 29.22757333333334),
    @(5, 8, 87.68272),
    @(5, 9, 0.08948272176993048),
    @(5, 10, 0.08948272176993047),
    @(5, 13, 2.688466333333334),
    @(5, 14, 8.065399000000001),
    @(5, 15, 0.1382584398787761),
    @(5, 16, 0.1382584398787761),
    @(5, 17, 78.5773469116978),
    @(5, 18, 707.1961222052802),
    @(5, 19, 0.01237174150801718),
    @(5, 20, 0.01237174150801718),
    @(6, 7, 29.22757333333334),
    @(6, 8, 87.68272),
    @(6, 9, 0.08948272176993048),
    @(6, 10, 0.08948272176993047),
    @(6, 13, 0.113399),
    @(6, 14, 0.340197),
    @(6, 15, 0.005831714769652435),
    @(6, 16, 0.005831714769652436),
    @(6, 17, 3.314377588426667),
    @(6, 18, 29.82939829584),
    @(6, 19, 0.000521837710174403),
    @(6, 20, 0.000521837710174403),
    @(7, 9, 0.4075568457508759),
    @(7, 10, 0.4075568457508759),
    @(7, 13, 15.47987166666667),
    @(7, 14, 46.439615),
    @(7, 15, 0.7960757698994193),
    @(7, 16, 0.7960757698994194),
    @(7, 17, 2060.673194287368),
    @(7, 18, 18546.05874858632),
    @(7, 19, 0.3244461297589074),
    @(7, 20, 0.3244461297589074),
    @(8, 9, 0.4075568457508759),
    @(8, 10, 0.4075568457508759),
    @(8, 15, 0.05100527512565552),
    @(8, 16, 0.05100527512565553),
    @(8, 19, 0.02078754904686778),
    @(8, 20, 0.02078754904686778),
    @(9, 9, 0.4075568457508759),
    @(9, 10, 0.4075568457508759),
    @(9, 13, 0.171678),
    @(9, 14, 0.515034),
    @(9, 15, 0.008828800326496623),
    @(9, 16, 0.008828800326496624),
    @(9, 17, 22.853694156306),
    @(9, 18, 205.683247406754),
    @(9, 19, 0.003598238012831267),
    @(9, 20, 0.003598238012831267),
    @(10, 9, 0.4075568457508759),
    @(10, 10, 0.4075568457508759),
    @(10, 13, 2.688466333333334),
    @(10, 14, 8.065399000000001),
    @(10, 15, 0.1382584398787761),
    @(10, 16, 0.1382584398787761),
    @(10, 17, 357.8873666487577),
    @(10, 18, 3220.986299838819),
    @(10, 19, 0.05634817365543109),
    @(10, 20, 0.05634817365543108),
    @(11, 9, 0.4075568457508759),
    @(11, 10, 0.4075568457508759),
    @(11, 13, 0.113399),
    @(11, 14, 0.340197),
    @(11, 15, 0.005831714769652435),
    @(11, 16, 0.005831714769652436),
    @(11, 17, 15.095621242273),
    @(11, 18, 135.860591180457),
    @(11, 19, 0.002376755276838342),
    @(11, 20, 0.002376755276838342),
    @(12, 7, 49.29039633333334),
    @(12, 8, 147.871189),
    @(12, 9, 0.1509067745968169),
    @(12, 10, 0.1509067745968168),
    @(12, 13, 15.47987166666667),
    @(12, 14, 46.439615),
    @(12, 15, 0.7960757698994193),
    @(12, 16, 0.7960757698994194),
    @(12, 17, 763.0090096391374),
    @(12, 18, 6867.081086752236),
    @(12, 19, 0.1201332267701991),
    @(12, 20, 0.1201332267701991),
    @(13, 7, 49.29039633333334),
    @(13, 8, 147.871189),
    @(13, 9, 0.1509067745968169),
    @(13, 10, 0.1509067745968168),
    @(13, 15, 0.05100527512565552),
    @(13, 16, 0.05100527512565553),
    @(13, 17, 48.886658696967),
    @(13, 18, 439.979928272703),
    @(13, 19, 0.007697041556635928),
    @(13, 20, 0.007697041556635927),
    @(14, 7, 49.29039633333334),
    @(14, 8, 147.871189),
    @(14, 9, 0.1509067745968169),
    @(14, 10, 0.1509067745968168),
    @(14, 13, 0.171678),
    @(14, 14, 0.515034),
    @(14, 15, 0.008828800326496623),
    @(14, 16, 0.008828800326496624),
    @(14, 17, 8.462076661714001),
    @(14, 18, 76.15868995542601),
    @(14, 19, 0.001332325780830929),
    @(14, 20, 0.001332325780830929),
    @(15, 7, 49.29039633333334),
    @(15, 8, 147.871189),
    @(15, 9, 0.1509067745968169),
    @(15, 10, 0.1509067745968168),
    @(15, 13, 2.688466333333334),
    @(15, 14, 8.065399000000001),
    @(15, 15, 0.1382584398787761),
    @(15, 16, 0.1382584398787761),
    @(15, 17, 132.5155710988235),
    @(15, 18, 1192.640139889411),
    @(15, 19, 0.02086413522289401),
    @(15, 20, 0.02086413522289401),
    @(16, 7, 49.29039633333334),
    @(16, 8, 147.871189),
    @(16, 9, 0.1509067745968169),
    @(16, 10, 0.1509067745968168),
    @(16, 13, 0.113399),
    @(16, 14, 0.340197),
    @(16, 15, 0.005831714769652435),
    @(16, 16, 0.005831714769652436),
    @(16, 17, 5.589481653803668),
    @(16, 18, 50.30533488423301),
    @(16, 19, 0.0008800452662568678),
    @(16, 20, 0.0008800452662568678),
    @(17, 7, 33.14535033333333),
    @(17, 8, 99.43605099999999),
    @(17, 9, 0.1014773319706963),
    @(17, 10, 0.1014773319706963),
    @(17, 13, 15.47987166666667),
    @(17, 14, 46.439615),
    @(17, 15, 0.7960757698994193),
    @(17, 16, 0.7960757698994194),
    @(17, 17, 513.0857695067073),
    @(17, 18, 4617.771925560365),
    @(17, 19, 0.08078364517591105),
    @(17, 20, 0.08078364517591105),
    @(18, 7, 33.14535033333333),
    @(18, 8, 99.43605099999999),
    @(18, 9, 0.1014773319706963),
    @(18, 10, 0.1014773319706963),
    @(18, 15, 0.05100527512565552),
    @(18, 16, 0.05100527512565553),
    @(18, 17, 32.873856768753),
    @(18, 18, 295.864710918777),
    @(18, 19, 0.005175879236182846),
    @(18, 20, 0.005175879236182846),
    @(19, 7, 33.14535033333333),
    @(19, 8, 99.43605099999999),
    @(19, 9, 0.1014773319706963),
    @(19, 10, 0.1014773319706963),
    @(19, 13, 0.171678),
    @(19, 14, 0.515034),
    @(19, 15, 0.008828800326496623),
    @(19, 16, 0.008828800326496624),
    @(19, 17, 5.690327454526),
    @(19, 18, 51.212947090734),
    @(19, 19, 0.00089592310163489),
    @(19, 20, 0.0008959231016348901),
    @(20, 7, 33.14535033333333),
    @(20, 8, 99.43605099999999),
    @(20, 9, 0.1014773319706963),
    @(20, 10, 0.1014773319706963),
    @(20, 13, 2.688466333333334),
    @(20, 14, 8.065399000000001),
    @(20, 15, 0.1382584398787761),
    @(20, 16, 0.1382584398787761),
    @(20, 17, 89.11015847770545),
    @(20, 18, 801.9914262993491),
    @(20, 19, 0.01403009760132912),
    @(20, 20, 0.01403009760132912),
    @(21, 7, 33.14535033333333),
    @(21, 8, 99.43605099999999),
    @(21, 9, 0.1014773319706963),
    @(21, 10, 0.1014773319706963),
    @(21, 13, 0.113399),
    @(21, 14, 0.340197),
    @(21, 15, 0.005831714769652435),
    @(21, 16, 0.005831714769652436),
    @(21, 17, 3.758649582449667),
    @(21, 18, 33.827846242047),
    @(21, 19, 0.0005917868556384331),
    @(21, 20, 0.0005917868556384331),
    @(22, 7, 81.84527466666667),
    @(22, 8, 245.535824),
    @(22, 9, 0.2505763259116804),
    @(22, 10, 0.2505763259116804),
    @(22, 13, 15.47987166666667),
    @(22, 14, 46.439615),
    @(22, 15, 0.7960757698994193),
    @(22, 16, 0.7960757698994194),
    @(22, 17, 1266.954348363085),
    @(22, 18, 11402.58913526776),
    @(22, 19, 0.1994777415687088),
    @(22, 20, 0.1994777415687088),
    @(23, 7, 81.84527466666667),
    @(23, 8, 245.535824),
    @(23, 9, 0.2505763259116804),
    @(23, 10, 0.2505763259116804),
    @(23, 15, 0.05100527512565552),
    @(23, 16, 0.05100527512565553),
    @(23, 17, 81.174880021872),
    @(23, 18, 730.5739201968479),
    @(23, 19, 0.01278071444310118),
    @(23, 20, 0.01278071444310118),
    @(24, 7, 81.84527466666667),
    @(24, 8, 245.535824),
    @(24, 9, 0.2505763259116804),
    @(24, 10, 0.2505763259116804),
    @(24, 13, 0.171678),
    @(24, 14, 0.515034),
    @(24, 15, 0.008828800326496623),
    @(24, 16, 0.008828800326496624),
    @(24, 17, 14.051033064224),
    @(24, 18, 126.459297578016),
    @(24, 19, 0.002212288348021368),
    @(24, 20, 0.002212288348021368),
    @(25, 7, 81.84527466666667),
    @(25, 8, 245.535824),
    @(25, 9, 0.2505763259116804),
    @(25, 10, 0.2505763259116804),
    @(25, 13, 2.688466333333334),
    @(25, 14, 8.065399000000001),
    @(25, 15, 0.1382584398787761),
    @(25, 16, 0.1382584398787761),
    @(25, 17, 220.0382654837529),
    @(25, 18, 1980.344389353776),
    @(25, 19, 0.03464429189110466),
    @(25, 20, 0.03464429189110466),
    @(26, 7, 81.84527466666667),
    @(26, 8, 245.535824),
    @(26, 9, 0.2505763259116804),
    @(26, 10, 0.2505763259116804),
    @(26, 13, 0.113399),
    @(26, 14, 0.340197),
    @(26, 15, 0.005831714769652435),
    @(26, 16, 0.005831714769652436),
    @(26, 17, 9.281172301925334),
    @(26, 18, 83.530550717328),
    @(26, 19, 0.001461289660744389),
    @(26, 20, 0.001461289660744389)
)

foreach ($u in $updates) {
    $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
}